$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '59.195.89'
$ws.Cells.Item(2, 5).Value = '  -2.07%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.582.88'
$ws.Cells.Item(3, 5).Value = '  -2.17%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.12%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''562.70'
$ws.Cells.Item(5, 5).Value = '  -1.21%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''142.68'
$ws.Cells.Item(6, 5).Value = '  -2.37%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.30%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -1.94%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.588.36'
$ws.Cells.Item(9, 5).Value = '  -2.83%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -2.98%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.96%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +11.60%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''0.349'
$ws.Cells.Item(13, 5).Value = '  +2.31%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '3.035.89'
$ws.Cells.Item(14, 5).Value = '  -2.52%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Avalanche'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(15, 4).Value = '''23.38'
$ws.Cells.Item(15, 5).Value = '  +7.26%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).Value = '59.151.14'
$ws.Cells.Item(16, 5).Value = '  -2.14%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.07%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.586.12'
$ws.Cells.Item(18, 5).Value = '  -2.66%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''4.57'
$ws.Cells.Item(19, 5).Value = '  +0.15%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''337.35'
$ws.Cells.Item(20, 5).Value = '  -2.11%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''10.36'
$ws.Cells.Item(21, 5).Value = '  -0.48%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''6.40'
$ws.Cells.Item(22, 5).Value = '  +0.24%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.20%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''64.09'
$ws.Cells.Item(24, 5).Value = '  -4.08%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''0.467'
$ws.Cells.Item(25, 5).Value = '  +5.83%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''0.999'
$ws.Cells.Item(26, 5).Value = '  +0.47%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -3.22%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''7.31'
$ws.Cells.Item(28, 5).Value = '  -0.50%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '0.0₃0773'
$ws.Cells.Item(29, 5).Value = '  -0.19%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.12%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(31, 4).Value = '''1.67'
$ws.Cells.Item(31, 5).Value = '  -2.69%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Monero'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(32, 4).Value = '''160.17'
$ws.Cells.Item(32, 5).Value = '  +2.53%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''6.10'
$ws.Cells.Item(33, 5).Value = '  -0.43%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''18.96'
$ws.Cells.Item(34, 5).Value = '  -1.30%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''4.03'
$ws.Cells.Item(35, 5).Value = '  -1.81%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -1.01%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -3.58%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -4.43%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.36%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -2.11%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''295.07'
$ws.Cells.Item(41, 5).Value = '  -3.58%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.17%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.998'
$ws.Cells.Item(43, 5).Value = '  +0.38%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''131.75'
$ws.Cells.Item(44, 5).Value = '  +3.90%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''0.0971'
$ws.Cells.Item(45, 5).Value = '  -1.13%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''0.595'
$ws.Cells.Item(46, 5).Value = '  -2.02%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Hedera'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(47, 4).Value = '''0.0536'
$ws.Cells.Item(47, 5).Value = '  -2.41%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(48, 4).Value = '''10.64'
$ws.Cells.Item(48, 5).Value = '  -0.15%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''19.00'
$ws.Cells.Item(49, 5).Value = '  -1.76%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -0.78%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''18.61'
$ws.Cells.Item(51, 5).Value = '  +0.54%  '
